$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row 6 label first: "Begründung"
$ws.Range("A6").Value = "Begründung"

# Update the "Lösungsvorschlag" text in B5 (wrap it, since it's now longer)
$ws.Range("B5").Value = "Passwörter für jeden Mitarbeiter einführen und komplexeres Admin-Passwort"
$ws.Range("B5").WrapText = $true

# Add the new row 6 justification text
$ws.Range("B6").Value = "Zur Sicherheit der Daten und eindeutigen Identifizierung durch Logs"
$ws.Range("B6").WrapText = $true

# Set row heights for the wrapped rows
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 30

# Update selection to match the new active cell
$ws.Range("B6").Select() | Out-Null
